# repull data, push all data, mean calculation
# Update column F (dSF) values for a set of rows to reflect re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 3
$ws.Range("F7").Value = -2
$ws.Range("F16").Value = 0
$ws.Range("F25").Value = -11
$ws.Range("F34").Value = -6
$ws.Range("F35").Value = -3
$ws.Range("F39").Value = 2
